$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 82  # H5
$ws.Cells.Item(5, 9).Value = 87.333336  # I5
$ws.Cells.Item(5, 10).Value = 50  # J5
$ws.Cells.Item(5, 11).Value = 87.333336  # K5
$ws.Cells.Item(5, 12).Value = 50  # L5
$ws.Cells.Item(5, 13).Value = 27.666664  # M5
$ws.Cells.Item(5, 14).Value = -280  # N5

$ws.Cells.Item(28, 8).Value = 1914.6666  # H28
$ws.Cells.Item(28, 9).Value = 1914.6666  # I28
$ws.Cells.Item(28, 10).Value = 0  # J28
$ws.Cells.Item(28, 11).Value = 1914.6666  # K28
$ws.Cells.Item(28, 12).Value = 0  # L28
$ws.Cells.Item(28, 13).Value = -1429.6666  # M28

$ws.Cells.Item(112, 8).Value = 3699.3333  # H112
$ws.Cells.Item(112, 9).Value = 1100  # I112
$ws.Cells.Item(112, 10).Value = 4999  # J112
$ws.Cells.Item(112, 11).Value = 3300  # K112
$ws.Cells.Item(112, 12).Value = 14997  # L112
$ws.Cells.Item(112, 13).Value = -2192  # M112
$ws.Cells.Item(112, 14).Value = -17213  # N112

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 0  # H41
$ws.Cells.Item(41, 9).Value = 0  # I41
$ws.Cells.Item(41, 10).Value = 0  # J41
$ws.Cells.Item(41, 11).Value = 0  # K41
$ws.Cells.Item(41, 12).Value = 0  # L41
$ws.Cells.Item(41, 13).ClearContents()  # M41

$ws.Cells.Item(132, 8).Value = 2374.9285  # H132
$ws.Cells.Item(132, 9).Value = 2270.8333  # I132
$ws.Cells.Item(132, 10).Value = 2999.5  # J132
$ws.Cells.Item(132, 11).Value = 6812.499899999999  # K132
$ws.Cells.Item(132, 12).Value = 8998.5  # L132
$ws.Cells.Item(132, 13).Value = -4282.499899999999  # M132
$ws.Cells.Item(132, 14).Value = -14058.5  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3192.0833  # H105
$ws.Cells.Item(105, 9).Value = 2970.5  # I105
$ws.Cells.Item(105, 10).Value = 4300  # J105
$ws.Cells.Item(105, 11).Value = 2970.5  # K105
$ws.Cells.Item(105, 12).Value = 4300  # L105
$ws.Cells.Item(105, 13).Value = -1223.5  # M105
$ws.Cells.Item(105, 14).Value = -7794  # N105

$ws.Cells.Item(134, 8).Value = 6575.5  # H134
$ws.Cells.Item(134, 9).Value = 3478.6667  # I134
$ws.Cells.Item(134, 10).Value = 9672.333000000001  # J134
$ws.Cells.Item(134, 11).Value = 10436.0001  # K134
$ws.Cells.Item(134, 12).Value = 29016.999  # L134
$ws.Cells.Item(134, 13).Value = -7901.000100000001  # M134
$ws.Cells.Item(134, 14).Value = -34086.999  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1351.6666  # H34
$ws.Cells.Item(34, 9).Value = 555  # I34
$ws.Cells.Item(34, 10).Value = 1750  # J34
$ws.Cells.Item(34, 11).Value = 1665  # K34
$ws.Cells.Item(34, 12).Value = 5250  # L34
$ws.Cells.Item(34, 13).Value = -1581  # M34
$ws.Cells.Item(34, 14).Value = -5418  # N34

$ws.Cells.Item(39, 8).Value = 4500  # H39
$ws.Cells.Item(39, 9).Value = 0  # I39
$ws.Cells.Item(39, 10).Value = 4500  # J39
$ws.Cells.Item(39, 11).Value = 0  # K39
$ws.Cells.Item(39, 12).Value = 13500  # L39
$ws.Cells.Item(39, 14).Value = -14088  # N39

$ws.Cells.Item(55, 8).Value = 1889.45  # H55
$ws.Cells.Item(55, 9).Value = 1378.9  # I55
$ws.Cells.Item(55, 10).Value = 2400  # J55
$ws.Cells.Item(55, 11).Value = 4136.700000000001  # K55
$ws.Cells.Item(55, 12).Value = 7200  # L55
$ws.Cells.Item(55, 13).Value = -3959.700000000001  # M55
$ws.Cells.Item(55, 14).Value = -7554  # N55

$ws.Cells.Item(131, 8).Value = 3602.8235  # H131
$ws.Cells.Item(131, 9).Value = 1800  # I131
$ws.Cells.Item(131, 10).Value = 4157.5386  # J131
$ws.Cells.Item(131, 11).Value = 5400  # K131
$ws.Cells.Item(131, 12).Value = 12472.6158  # L131
$ws.Cells.Item(131, 13).Value = -360  # M131
$ws.Cells.Item(131, 14).Value = -22552.6158  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4619.7144  # H102
$ws.Cells.Item(102, 9).Value = 5662.6665  # I102
$ws.Cells.Item(102, 10).Value = 3837.5  # J102
$ws.Cells.Item(102, 11).Value = 5662.6665  # K102
$ws.Cells.Item(102, 12).Value = 3837.5  # L102
$ws.Cells.Item(102, 13).Value = -4040.6665  # M102
$ws.Cells.Item(102, 14).Value = -7081.5  # N102

$ws.Cells.Item(113, 8).Value = 500  # H113
$ws.Cells.Item(113, 9).Value = 500  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 500  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).Value = 1670  # M113

$ws.Cells.Item(125, 8).Value = 0  # H125
$ws.Cells.Item(125, 9).Value = 0  # I125
$ws.Cells.Item(125, 10).Value = 0  # J125
$ws.Cells.Item(125, 11).Value = 0  # K125
$ws.Cells.Item(125, 12).Value = 0  # L125

$ws.Cells.Item(126, 8).Value = 5666.6665  # H126
$ws.Cells.Item(126, 9).Value = 4000  # I126
$ws.Cells.Item(126, 10).Value = 9000  # J126
$ws.Cells.Item(126, 11).Value = 12000  # K126
$ws.Cells.Item(126, 12).Value = 27000  # L126
$ws.Cells.Item(126, 13).Value = -9530  # M126
$ws.Cells.Item(126, 14).Value = -31940  # N126

$ws.Cells.Item(127, 8).Value = 0  # H127
$ws.Cells.Item(127, 9).Value = 0  # I127
$ws.Cells.Item(127, 10).Value = 0  # J127
$ws.Cells.Item(127, 11).Value = 0  # K127
$ws.Cells.Item(127, 12).Value = 0  # L127

$ws.Cells.Item(128, 8).Value = 0  # H128
$ws.Cells.Item(128, 9).Value = 0  # I128
$ws.Cells.Item(128, 10).Value = 0  # J128
$ws.Cells.Item(128, 11).Value = 0  # K128
$ws.Cells.Item(128, 12).Value = 0  # L128

$ws.Cells.Item(129, 8).Value = 0  # H129
$ws.Cells.Item(129, 9).Value = 0  # I129
$ws.Cells.Item(129, 10).Value = 0  # J129
$ws.Cells.Item(129, 11).Value = 0  # K129
$ws.Cells.Item(129, 12).Value = 0  # L129

$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 9).Value = 0  # I130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 11).Value = 0  # K130
$ws.Cells.Item(130, 12).Value = 0  # L130

$ws.Cells.Item(131, 8).Value = 0  # H131
$ws.Cells.Item(131, 9).Value = 0  # I131
$ws.Cells.Item(131, 10).Value = 0  # J131
$ws.Cells.Item(131, 11).Value = 0  # K131
$ws.Cells.Item(131, 12).Value = 0  # L131

$ws.Cells.Item(132, 8).Value = 3346  # H132
$ws.Cells.Item(132, 9).Value = 2756.4  # I132
$ws.Cells.Item(132, 10).Value = 5311.3335  # J132
$ws.Cells.Item(132, 11).Value = 8269.200000000001  # K132
$ws.Cells.Item(132, 12).Value = 15934.0005  # L132
$ws.Cells.Item(132, 13).Value = -5739.200000000001  # M132
$ws.Cells.Item(132, 14).Value = -20994.0005  # N132

$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 9).Value = 0  # I133
$ws.Cells.Item(133, 10).Value = 0  # J133
$ws.Cells.Item(133, 11).Value = 0  # K133
$ws.Cells.Item(133, 12).Value = 0  # L133

$ws.Cells.Item(134, 8).Value = 0  # H134
$ws.Cells.Item(134, 9).Value = 0  # I134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 11).Value = 0  # K134
$ws.Cells.Item(134, 12).Value = 0  # L134

$ws.Cells.Item(135, 8).Value = 50000  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 50000  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 50000  # L135
$ws.Cells.Item(135, 14).Value = -60140  # N135

$ws.Cells.Item(136, 8).Value = 40000  # H136
$ws.Cells.Item(136, 9).Value = 0  # I136
$ws.Cells.Item(136, 10).Value = 40000  # J136
$ws.Cells.Item(136, 11).Value = 0  # K136
$ws.Cells.Item(136, 12).Value = 120000  # L136
$ws.Cells.Item(136, 14).Value = -125100  # N136

$ws.Cells.Item(137, 8).Value = 0  # H137
$ws.Cells.Item(137, 9).Value = 0  # I137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 11).Value = 0  # K137
$ws.Cells.Item(137, 12).Value = 0  # L137

$ws.Cells.Item(138, 8).Value = 200000  # H138
$ws.Cells.Item(138, 9).Value = 200000  # I138
$ws.Cells.Item(138, 10).Value = 0  # J138
$ws.Cells.Item(138, 11).Value = 200000  # K138
$ws.Cells.Item(138, 12).Value = 0  # L138
$ws.Cells.Item(138, 13).Value = -194860  # M138

$ws.Cells.Item(139, 8).Value = 0  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 0  # L139

$ws.Cells.Item(140, 8).Value = 150000  # H140
$ws.Cells.Item(140, 9).Value = 0  # I140
$ws.Cells.Item(140, 10).Value = 150000  # J140
$ws.Cells.Item(140, 11).Value = 0  # K140
$ws.Cells.Item(140, 12).Value = 150000  # L140
$ws.Cells.Item(140, 14).Value = -160360  # N140

$ws.Cells.Item(141, 8).Value = 0  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 0  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 0  # L141

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 4619.619  # H100
$ws.Cells.Item(100, 9).Value = 2356.125  # I100
$ws.Cells.Item(100, 10).Value = 11862.8  # J100
$ws.Cells.Item(100, 11).Value = 2356.125  # K100
$ws.Cells.Item(100, 12).Value = 11862.8  # L100
$ws.Cells.Item(100, 13).Value = -1815.125  # M100
$ws.Cells.Item(100, 14).Value = -12944.8  # N100

$ws.Cells.Item(122, 8).Value = 5600.8  # H122
$ws.Cells.Item(122, 9).Value = 4668  # I122
$ws.Cells.Item(122, 10).Value = 7000  # J122
$ws.Cells.Item(122, 11).Value = 14004  # K122
$ws.Cells.Item(122, 12).Value = 21000  # L122
$ws.Cells.Item(122, 13).Value = -11554  # M122
$ws.Cells.Item(122, 14).Value = -25900  # N122

$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 9).Value = 0  # I124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 11).Value = 0  # K124
$ws.Cells.Item(124, 12).Value = 0  # L124

$ws.Cells.Item(125, 8).Value = 0  # H125
$ws.Cells.Item(125, 9).Value = 0  # I125
$ws.Cells.Item(125, 10).Value = 0  # J125
$ws.Cells.Item(125, 11).Value = 0  # K125
$ws.Cells.Item(125, 12).Value = 0  # L125

$ws.Cells.Item(127, 8).Value = 79999  # H127
$ws.Cells.Item(127, 9).Value = 0  # I127
$ws.Cells.Item(127, 10).Value = 79999  # J127
$ws.Cells.Item(127, 11).Value = 0  # K127
$ws.Cells.Item(127, 12).Value = 79999  # L127
$ws.Cells.Item(127, 14).Value = -89919  # N127

$ws.Cells.Item(128, 8).Value = 0  # H128
$ws.Cells.Item(128, 9).Value = 0  # I128
$ws.Cells.Item(128, 10).Value = 0  # J128
$ws.Cells.Item(128, 11).Value = 0  # K128
$ws.Cells.Item(128, 12).Value = 0  # L128

$ws.Cells.Item(129, 8).Value = 0  # H129
$ws.Cells.Item(129, 9).Value = 0  # I129
$ws.Cells.Item(129, 10).Value = 0  # J129
$ws.Cells.Item(129, 11).Value = 0  # K129
$ws.Cells.Item(129, 12).Value = 0  # L129

$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 9).Value = 0  # I130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 11).Value = 0  # K130
$ws.Cells.Item(130, 12).Value = 0  # L130

$ws.Cells.Item(131, 8).Value = 0  # H131
$ws.Cells.Item(131, 9).Value = 0  # I131
$ws.Cells.Item(131, 10).Value = 0  # J131
$ws.Cells.Item(131, 11).Value = 0  # K131
$ws.Cells.Item(131, 12).Value = 0  # L131

$ws.Cells.Item(132, 8).Value = 3579.8333  # H132
$ws.Cells.Item(132, 9).Value = 4698  # I132
$ws.Cells.Item(132, 10).Value = 1343.5  # J132
$ws.Cells.Item(132, 11).Value = 14094  # K132
$ws.Cells.Item(132, 12).Value = 4030.5  # L132
$ws.Cells.Item(132, 13).Value = -11564  # M132
$ws.Cells.Item(132, 14).Value = -9090.5  # N132

$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 9).Value = 0  # I133
$ws.Cells.Item(133, 10).Value = 0  # J133
$ws.Cells.Item(133, 11).Value = 0  # K133
$ws.Cells.Item(133, 12).Value = 0  # L133

$ws.Cells.Item(134, 8).Value = 0  # H134
$ws.Cells.Item(134, 9).Value = 0  # I134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 11).Value = 0  # K134
$ws.Cells.Item(134, 12).Value = 0  # L134

$ws.Cells.Item(135, 8).Value = 0  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 0  # L135

$ws.Cells.Item(136, 8).Value = 3995.6667  # H136
$ws.Cells.Item(136, 9).Value = 3993.5  # I136
$ws.Cells.Item(136, 10).Value = 4000  # J136
$ws.Cells.Item(136, 11).Value = 11980.5  # K136
$ws.Cells.Item(136, 12).Value = 12000  # L136
$ws.Cells.Item(136, 13).Value = -9430.5  # M136
$ws.Cells.Item(136, 14).Value = -17100  # N136

$ws.Cells.Item(137, 8).Value = 55000  # H137
$ws.Cells.Item(137, 9).Value = 0  # I137
$ws.Cells.Item(137, 10).Value = 55000  # J137
$ws.Cells.Item(137, 11).Value = 0  # K137
$ws.Cells.Item(137, 12).Value = 55000  # L137
$ws.Cells.Item(137, 14).Value = -65200  # N137

$ws.Cells.Item(138, 8).Value = 0  # H138
$ws.Cells.Item(138, 9).Value = 0  # I138
$ws.Cells.Item(138, 10).Value = 0  # J138
$ws.Cells.Item(138, 11).Value = 0  # K138
$ws.Cells.Item(138, 12).Value = 0  # L138

$ws.Cells.Item(139, 8).Value = 0  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 0  # L139

$ws.Cells.Item(140, 8).Value = 80000  # H140
$ws.Cells.Item(140, 9).Value = 0  # I140
$ws.Cells.Item(140, 10).Value = 80000  # J140
$ws.Cells.Item(140, 11).Value = 0  # K140
$ws.Cells.Item(140, 12).Value = 80000  # L140
$ws.Cells.Item(140, 14).Value = -90360  # N140

$ws.Cells.Item(141, 8).Value = 200000  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 200000  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 200000  # L141
$ws.Cells.Item(141, 14).Value = -210360  # N141
